# Daily attendance processing - 2026-01-07 11:06:19
# Rotate the "Recorded By" (column G) list of names left by one position
# for every row whose value contains more than one comma-separated name,
# except rows whose value already equals "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1 -and $text -ne "dnasr281@gmail.com, System") {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
    }
}
